# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) updates for column F
$updates = @{
    "展览" = @(
        @{ Row = 2;  Value = 240 },
        @{ Row = 3;  Value = 1088 },
        @{ Row = 5;  Value = 423 },
        @{ Row = 7;  Value = 549 },
        @{ Row = 9;  Value = 6759 },
        @{ Row = 10; Value = 155 },
        @{ Row = 11; Value = 94 },
        @{ Row = 15; Value = 1086 },
        @{ Row = 16; Value = 16137 },
        @{ Row = 17; Value = 1581 },
        @{ Row = 18; Value = 37 },
        @{ Row = 20; Value = 177 },
        @{ Row = 22; Value = 11319 },
        @{ Row = 24; Value = 945 },
        @{ Row = 25; Value = 4449 },
        @{ Row = 26; Value = 305 },
        @{ Row = 29; Value = 41 }
    )
    "全部类型" = @(
        @{ Row = 2;  Value = 240 },
        @{ Row = 3;  Value = 1088 },
        @{ Row = 5;  Value = 423 },
        @{ Row = 7;  Value = 549 },
        @{ Row = 10; Value = 6759 },
        @{ Row = 11; Value = 155 },
        @{ Row = 12; Value = 94 },
        @{ Row = 17; Value = 1086 },
        @{ Row = 18; Value = 16137 },
        @{ Row = 19; Value = 1581 },
        @{ Row = 20; Value = 37 },
        @{ Row = 22; Value = 177 },
        @{ Row = 26; Value = 11319 },
        @{ Row = 28; Value = 945 },
        @{ Row = 29; Value = 4449 },
        @{ Row = 30; Value = 305 },
        @{ Row = 33; Value = 41 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Cells.Item($u.Row, 6).Value = $u.Value
    }
}
